# "fix: ordine delle immagini nella presentazione"
# Rearranges the four pictures on slide 6 (the "ISTANZE EUCLIDEE" slide)
# into their corrected 2x2 grid positions.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# "Immagine 8" (top-left) -> moves to bottom-left
$sh = $s.Shapes.Item("Immagine 8")
$sh.Left = 0.00
$sh.Top = 290.5964

# "Segnaposto contenuto 4" (top-right) -> slides to top-left
$sh = $s.Shapes.Item("Segnaposto contenuto 4")
$sh.Left = 0.79402

# "Immagine 6" (bottom-right) -> moves up to top-right
$sh = $s.Shapes.Item("Immagine 6")
$sh.Left = 470.7099
$sh.Top = 72.2558

# "Immagine 10" (bottom-left) -> slides to bottom-right
$sh = $s.Shapes.Item("Immagine 10")
$sh.Left = 469.3219
